$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $c = $ws.Range($range)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.798.54"
$ws.Range("E2").Value = "  -2.43%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.568.49"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "206.67"
$ws.Range("E5").Value = "  -0.96%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -2.27%  "

# Row 8 - Solana
$ws.Range("E8").Value = "  -0.79%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.75%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.25%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.23%  "

# Row 12 - Wrapped liquid staked Ether 2.0
Set-TextValue "D12" "1.791.81"
$ws.Range("E12").Value = "  +0.01%  "

# Row 13 - Wrapped Ether
Set-TextValue "D13" "1.561.41"
$ws.Range("E13").Value = "  -0.55%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.42%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.82%  "

# Row 16 - Wrapped BTC
Set-TextValue "D16" "26.799.15"
$ws.Range("E16").Value = "  -2.43%  "

# Row 17 - Litecoin
Set-TextValue "D17" "61.43"
$ws.Range("E17").Value = "  -3.70%  "

# Row 18 - Chainlink
Set-TextValue "D18" "7.43"
$ws.Range("E18").Value = "  +1.88%  "

# Row 19 - Bitcoin Cash
Set-TextValue "D19" "215.15"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20 - Shiba Inu
$ws.Range("D20").Value = "0.0₃0678"
$ws.Range("E20").Value = "  -2.04%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.02%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  -2.68%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.87%  "

# Row 25 - Monero
Set-TextValue "D25" "153.02"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +0.54%  "

# Row 27 - Ethereum Classic
$ws.Range("E27").Value = "  -0.37%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -1.41%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.12"
$ws.Range("E31").Value = "  -3.14%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.22%  "

# Row 33 - Maker
Set-TextValue "D33" "1.396.83"
$ws.Range("E33").Value = "  +1.29%  "

# Row 34 - Internet Computer (DFINITY)
$ws.Range("E34").Value = "  -1.23%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -1.02%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.25%  "

# Row 37 - TrustWalletToken
$ws.Range("E37").Value = "  -1.95%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -2.93%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -2.70%  "

# Row 40 - ARBITRUM
Set-TextValue "D40" "0.816"
$ws.Range("E40").Value = "  -1.36%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - WEMIXToken
Set-TextValue "D42" "0.987"
$ws.Range("E42").Value = "  +0.82%  "

# Row 43 - RenderToken
Set-TextValue "D43" "1.81"
$ws.Range("E43").Value = "  +0.03%  "

# Row 44 - swaps from FraxShare to MXToken
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D44" "2.19"
$ws.Range("E44").Value = "  +1.10%  "

# Row 45 - swaps from MXToken to FraxShare
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D45" "5.33"
$ws.Range("E45").Value = "  +1.04%  "

# Row 46 - Aave
Set-TextValue "D46" "63.37"

# Row 47 - RocketPoolETH
Set-TextValue "D47" "1.704.23"
$ws.Range("E47").Value = "  +0.17%  "

# Row 48 - Quant
$ws.Range("E48").Value = "  +0.77%  "

# Row 49 - BabyDogeCoin
$ws.Range("D49").Value = "0.0₇0985"

# Row 50 - Algorand
$ws.Range("E50").Value = "  -0.31%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -0.88%  "
